$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$total = 0
for ($row = 3; $row -le 33; $row++) {
    $kVal = $ws.Cells.Item($row, 11).Value2
    $ws.Cells.Item($row, 12).Value = $kVal
    $total = $total + $kVal
}

$ws.Range("C47").Value = $total
